# Apply the "TODO: finish blast shield and potentially flip servo mount" update
# to the tasks tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Redesign avionics bay): update the comment text ---
# Old: "Everything should fit, ... Also unsure if I can stack MOSFET boards ... INCLUDE SERVO (forgot)"
# New: "Everything should fit, ... Accidentally modeled servo retention on wrong side of tube I think, gotta flip it"
$ws.Range("D4").Value2 = "Everything should fit, need to double check once I get components. Accidentally modeled servo retention on wrong side of tube I think, gotta flip it"

# --- Row 5 (Design and install blast shield): NOT DONE -> IN PROGRESS, add comment ---
# Copy the "IN PROGRESS" look (style + date number format) from row 4, which is already IN PROGRESS.
$ws.Range("A5").Style = $ws.Range("A4").Style

$ws.Range("B5").Style = $ws.Range("B4").Style
$ws.Range("B5").NumberFormat = $ws.Range("B4").NumberFormat

$ws.Range("C5").Style = $ws.Range("C4").Style
$ws.Range("C5").Value2 = $ws.Range("C4").Value2

$ws.Range("D5").Style = $ws.Range("D4").Style
$ws.Range("D5").Value2 = "Blast shield redesigned to two piece unit with wire passthrough. Need to CAD screwholes and then should be good to install"

# --- Row 10 (Figure out main parachute retention system): NOT DONE -> IN PROGRESS, add comment ---
$ws.Range("A10").Style = $ws.Range("A4").Style

$ws.Range("B10").Style = $ws.Range("B4").Style
$ws.Range("B10").NumberFormat = $ws.Range("B4").NumberFormat

$ws.Range("C10").Style = $ws.Range("C4").Style
$ws.Range("C10").Value2 = $ws.Range("C4").Value2

$ws.Range("D10").Style = $ws.Range("D4").Style
$ws.Range("D10").Value2 = "Prototype designed and retaining servo integrated in avionics bay. Testing tbd"

# --- Update the active selection/view to rest on D10, as in the saved workbook ---
$ws.Range("D10").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
